{"js": "// \"download articles with pandoc title blocks\"\n//\n// The article used to open with a bookmarked Heading-1 title paragraph\n// (\"Another Letter to an Agnostic\") followed by a bold \"By Dorothy Day\"\n// byline paragraph. Pandoc-style title blocks instead want a single\n// Title-styled paragraph holding just the author's name, so:\n//   1. the old heading paragraph (and its bookmark) goes away entirely\n//   2. the byline paragraph becomes a \"Title\" styled paragraph whose\n//      text is split into \"Dorothy\" / \" \" / \"Day\" runs (no more bold,\n//      no more \"By \" prefix).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Best-effort: the original markup also wraps the old title paragraph in\n// <w:bookmarkStart/.../w:bookmarkEnd name=\"another-letter-to-an-agnostic\">.\n// Try to drop it through the supported bookmark API before touching the\n// paragraphs themselves (harmless no-op on hosts that don't expose it).\ntry {\n  doc.deleteBookmark(\"another-letter-to-an-agnostic\");\n} catch (e) {\n  // ignore - not every host implements bookmark deletion\n}\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Paragraph 0 is the \"Another Letter to an Agnostic\" Heading1 paragraph -\n// remove it completely (text + paragraph mark).\nparagraphs.items[0].delete();\nawait context.sync();\n\n// Re-fetch: the former paragraph 1 (\"By Dorothy Day\") is now first.\nconst bylineParagraphs = body.paragraphs;\nbylineParagraphs.load(\"text\");\nawait context.sync();\n\nconst byline = bylineParagraphs.items[0];\n\n// Replace the whole paragraph (its bold \"By Dorothy Day\" run included)\n// with a Title-styled paragraph made of three plain runs: \"Dorothy\",\n// \" \" and \"Day\" - matching the target markup exactly.\nconst titleOoxml =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Dorothy</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">Day</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nbyline.getRange(\"Whole\").insertOoxml(titleOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"download articles with pandoc title blocks\"\n#\n# The article used to open with a bookmarked Heading-1 title paragraph\n# (\"Another Letter to an Agnostic\") followed by a bold \"By Dorothy Day\"\n# byline paragraph. Pandoc-style title blocks instead want a single\n# Title-styled paragraph holding just the author's name, so:\n#   1. the old heading paragraph (and its bookmark) goes away entirely\n#   2. the byline paragraph becomes a \"Title\" styled paragraph whose\n#      text is split into \"Dorothy\" / \" \" / \"Day\" runs (no more bold,\n#      no more \"By \" prefix).\n\n$d = $word.ActiveDocument\n\n# Best-effort: the original markup also wraps the old title paragraph in\n# <w:bookmarkStart/.../w:bookmarkEnd name=\"another-letter-to-an-agnostic\">.\n# Try to drop it through the supported Bookmarks API before touching the\n# paragraphs themselves (harmless no-op on hosts that don't implement it).\ntry {\n    $d.Bookmarks.Item(\"another-letter-to-an-agnostic\").Delete()\n} catch {\n    # ignore - not every host implements bookmark deletion\n}\n\n# Paragraph 1 is the \"Another Letter to an Agnostic\" Heading1 paragraph -\n# remove it completely (text + paragraph mark).\n$titleParagraph = $d.Paragraphs.Item(1)\n$titleParagraph.Range.Delete()\n\n# The former paragraph 2 (\"By Dorothy Day\") is now first.\n$byline = $d.Paragraphs.Item(1)\n\n# Replace the whole paragraph (its bold \"By Dorothy Day\" run included)\n# with a Title-styled paragraph made of three plain runs: \"Dorothy\",\n# \" \" and \"Day\" - matching the target markup exactly.\n$titleXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr><w:r><w:t xml:space=\"preserve\">Dorothy</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$byline.Range.InsertXML($titleXml)\n"}
